$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly record was added for this market/product, inserted right
# before the former row 198 (everything below shifts down by one row,
# and the sheet's used range grows from R276 to R277).
$ws.Rows.Item(198).Insert()

$ws.Range("A198").Value = 4
$ws.Range("B198").Value = 'Feria Lagunitas de Puerto Montt'
$ws.Range("C198").Value = 'Los Lagos'
$ws.Range("D198").Value = 44784
$ws.Range("E198").Value = 10
$ws.Range("F198").Value = 100112032
$ws.Range("G198").Value = 'Zapallo italiano'
$ws.Range("H198").Value = 'Sin especificar'
$ws.Range("I198").Value = 'Primera'
$ws.Range("J198").Value = 70
$ws.Range("K198").Value = 26000
$ws.Range("L198").Value = 26000
$ws.Range("M198").Value = 26000
$ws.Range("N198").Value = '$/caja 50 unidades'
$ws.Range("O198").Value = 'Región de Arica y Parinacota'
$ws.Range("P198").Value = 520
$ws.Range("Q198").Value = 50
$ws.Range("R198").Value = 'Hortaliza'
